$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width change (46.33203125 -> ~61.33203125) ---
$ws.Columns.Item(2).ColumnWidth = 60.5

# --- Row 46 (the "test" row): rename to "test entry", fill in repository (C46) ---
# and publication (D46) with the same link as website, move language value to "test"
$ws.Cells.Item(46, 1).Value = "test entry"
$ws.Cells.Item(46, 3).Value = $ws.Cells.Item(46, 2).Value()
$ws.Cells.Item(46, 3).HorizontalAlignment = -4131
$ws.Cells.Item(46, 4).Value = $ws.Cells.Item(46, 2).Value()
$ws.Cells.Item(46, 4).HorizontalAlignment = -4131
$ws.Cells.Item(46, 5).Value = "test"

# --- Row 47 (new row): simplerspec entry ---
$ws.Cells.Item(47, 1).Value = "simplerspec"
$ws.Cells.Item(47, 3).Value = "https://github.com/philipp-baumann/simplerspec"
$ws.Cells.Item(47, 2).Value = "https://github.com/philipp-baumann/simplerspec-pedometron-article"
$ws.Cells.Item(47, 5).Value = "R"

# --- Row 45 ("jcamp"): clear website (B45), add repository link (C45) ---
$ws.Cells.Item(45, 2).ClearContents()
$ws.Cells.Item(45, 3).Value = "https://github.com/nzhagen/jcamp"

# --- Selection / scroll position ---
$ws.Range("E46").Select()
